$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 14
$ws.Range("E4").Value = 7
$ws.Range("E20").Value = 26
$ws.Range("E21").Value = 27
$ws.Range("E26").Value = 12
$ws.Range("E27").Value = 42
$ws.Range("E29").Value = 10
$ws.Range("E31").Value = 4
$ws.Range("E33").Value = 11
